$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: introduce brand-new shared strings in the exact order they must
# appear in xl/sharedStrings.xml (the engine appends newly-seen strings to
# the shared string table in first-write order).
# ---------------------------------------------------------------------------
$ws.Cells.Item(5,4).Value  = "After copy contructor and"   # D5  -> new string
$ws.Cells.Item(9,4).Value  = "using member function"       # D9  -> new string
$ws.Cells.Item(3,4).Value  = "no=10"                        # D3  -> new string
$ws.Cells.Item(7,4).Value  = "no=11"                        # D7  -> new string
$ws.Cells.Item(2,3).Value  = "no=9"                         # C2  -> new string
$ws.Cells.Item(3,3).Value  = "name=raj"                     # C3  -> new string
$ws.Cells.Item(14,2).Value = "to display help command"      # B14 -> new string
$ws.Cells.Item(14,3).Value = "argv[0] -h"                    # C14 -> new string
$ws.Cells.Item(14,4).Value = "enter given inputs"            # D14 -> new string
$ws.Cells.Item(17,2).Value = "to check with char"            # B17 -> new string
$ws.Cells.Item(17,3).Value = "no=a"                          # C17 -> new string
$ws.Cells.Item(24,4).Value = "name=garbage"                  # D24 -> new string
$ws.Cells.Item(25,4).Value = "exit of program"               # D25 -> new string

# ---------------------------------------------------------------------------
# Step 2: fill in the rest of the grid (numbers + reused strings). Order no
# longer matters for shared-string positions because every unique string has
# already been registered above.
# ---------------------------------------------------------------------------

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,4).Value = "using parameterized constructor"
$ws.Cells.Item(2,5).Value = "using parameterized constructor"

# --- Row 3 ---
$ws.Cells.Item(3,5).Value = "no=10"

# --- Row 4 ---
$ws.Cells.Item(4,4).Value = "name=prasuna"
$ws.Cells.Item(4,5).Value = "name=prasuna"
$ws.Cells.Item(4,6).Value = "PASS"

# --- Row 5 ---
$ws.Cells.Item(5,5).Value = "After copy contructor and"

# --- Row 6 ---
$ws.Cells.Item(6,4).Value = "after operator overloading"
$ws.Cells.Item(6,5).Value = "after operator overloading"

# --- Row 7 ---
$ws.Cells.Item(7,5).Value = "no=11"

# --- Row 8 (new row) ---
$ws.Cells.Item(8,4).Value = "name=prasuna"
$ws.Cells.Item(8,5).Value = "name=prasuna"

# --- Row 9 ---
$ws.Cells.Item(9,5).Value = "using member function"

# --- Row 10 (new row) ---
$ws.Cells.Item(10,4).Value = "no=9"
$ws.Cells.Item(10,5).Value = "no=9"

# --- Row 11 (new row) ---
$ws.Cells.Item(11,4).Value = "name=raj"
$ws.Cells.Item(11,5).Value = "name=raj"

# --- Row 14 (new test case 2) ---
$ws.Cells.Item(14,1).Value = 2
$ws.Cells.Item(14,5).Value = "enter given inputs"
$ws.Cells.Item(14,6).Value = "PASS"

# --- Row 17 (new test case 3) ---
$ws.Cells.Item(17,1).Value = 3
$ws.Cells.Item(17,4).Value = "using parameterized constructor"
$ws.Cells.Item(17,5).Value = "using parameterized constructor"

# --- Row 18 (new row) ---
$ws.Cells.Item(18,4).Value = "no=10"
$ws.Cells.Item(18,5).Value = "no=10"

# --- Row 19 (new row) ---
$ws.Cells.Item(19,4).Value = "name=prasuna"
$ws.Cells.Item(19,5).Value = "name=prasuna"

# --- Row 20 (new row) ---
$ws.Cells.Item(20,4).Value = "After copy contructor and"
$ws.Cells.Item(20,5).Value = "After copy contructor and"
$ws.Cells.Item(20,6).Value = "PASS"

# --- Row 21 (new row) ---
$ws.Cells.Item(21,4).Value = "after operator overloading"
$ws.Cells.Item(21,5).Value = "after operator overloading"

# --- Row 22 (new row) ---
$ws.Cells.Item(22,4).Value = "no=11"
$ws.Cells.Item(22,5).Value = "no=11"

# --- Row 23 (new row) ---
$ws.Cells.Item(23,4).Value = "name=prasuna"
$ws.Cells.Item(23,5).Value = "name=prasuna"

# --- Row 24 ---
$ws.Cells.Item(24,5).Value = "name=garbage"

# --- Row 25 ---
$ws.Cells.Item(25,5).Value = "exit of program"

# ---------------------------------------------------------------------------
# Step 3: update the view state (scroll position + selection) to match the
# saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("F21").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
